$wb = $excel.ActiveWorkbook

# --- Insert the new "CoreBankRegularizationSummary" sheet ---
# Positioned right after "CoreBankAllocationSummary" (i.e. right before
# "Updation_of_Disposition"), matching the new sheet order / sheetId / rIds.
$src = $wb.Worksheets.Item("CoreBankAllocationSummary")
$newWs = $wb.Worksheets.Add($wb.Worksheets.Item("Updation_of_Disposition"))
$newWs.Name = "CoreBankRegularizationSummary"

# Copy header + first data row (with formatting) from CoreBankAllocationSummary:
#   A1:B2 (TestScenario/Run columns) -> A1:B2
#   E1:F2 (Region/Branch columns)    -> C1:D2
$src.Range("A1:B2").Copy($newWs.Range("A1:B2"))
$src.Range("E1:F2").Copy($newWs.Range("C1:D2"))

# This new sheet's own TestScenario name replaces the copied one
$newWs.Range("A2").Value = "CoreBankRegularizationSummary"

# Column widths (approximate best-fit match; engine quantizes to 1/6 steps)
$newWs.Range("A1").EntireColumn.ColumnWidth = 30.33
$newWs.Range("C1").EntireColumn.ColumnWidth = 8
$newWs.Range("D1").EntireColumn.ColumnWidth = 8.5

# --- CoreBankAllocationSummary loses the active-tab/selection it had ---
$src.Range("F13").Select()

# --- New sheet becomes the active / selected tab ---
$newWs.Range("F17").Select()
